$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement table (rows 16-22) is being refreshed: previous
# worker/period entries are swapped out for the new data, and a couple of
# "Valor Mora" / "Salario Basico" amounts are corrected.

# Row 16: was HERNAN DAVID MONTERROZA PIMENTEL / 2208 -> now YUDIS DEL CARMEN ALCALA GONZALEZ / 2303
$ws.Range("C16").Value = "30873761"
$ws.Range("D16").Value = "YUDIS DEL CARMEN ALCALA GONZALEZ"
$ws.Range("E16").Value = "2303"
$ws.Range("F16").Value = 46400
$ws.Range("G16").Value = 1423500

# Row 19: was YUDIS DEL CARMEN ALCALA GONZALEZ / 2303 -> now HERNAN DAVID MONTERROZA PIMENTEL / 2208
$ws.Range("C19").Value = "1002191794"
$ws.Range("D19").Value = "HERNAN DAVID MONTERROZA PIMENTEL"
$ws.Range("E19").Value = "2208"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1300000

# Row 20: period 2304 -> 2305, Valor Mora 30934 -> 17013
$ws.Range("E20").Value = "2305"
$ws.Range("F20").Value = 17013

# Row 21: period 2305 -> 2304, Valor Mora 17013 -> 30934
$ws.Range("E21").Value = "2304"
$ws.Range("F21").Value = 30934
